$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the formula error: columns C (price) and F (real multiplier) were
# hard-coded values; replace them with proper formulas derived from the
# exchange rate (D) and official-group multiplier (E) columns.
# C = 25 * D * E
# F = C / 25
# Row 7 (NekoCode) also had its exchange-rate cell (D7) corrected from the
# wrong hard-coded "1" to the proper 0.5.

$ws.Range("D7").Value = 0.5

for ($r = 3; $r -le 13; $r++) {
    $ws.Range("C$r").Formula = "=25*D$r*E$r"
    $ws.Range("F$r").Formula = "=C$r/25"
}

# Excel auto-applied a 2-decimal number format to C13 after editing it.
$ws.Range("C13").NumberFormat = "0.00"

# Move the active selection to C13 (previously C14)
[void]$ws.Range("C13").Select()
